# Update countries & provincias Spain
# Applies the 7-Aug-2020 data refresh to the "Pais" sheet:
#   - bumps the "Datos actualizados..." timestamp to 12:42
#   - refreshes case/recovered/death counts for a handful of countries
#   - Malta overtakes Jamaica (row 151/152) and Timor Oriental overtakes
#     Santa Lucia (row 202/203) in the ranking, so those two row pairs swap
#     which country name they display

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 7 de Agosto de 2020 a las 12:42"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5032805
$ws.Range("C4").Value = 626
$ws.Range("E4").Value = 2292644
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 162812

# Iran (row 14)
$ws.Range("B14").Value = 322567
$ws.Range("C14").Value = 2450
$ws.Range("D14").Value = 279724
$ws.Range("E14").Value = 24711
$ws.Range("G14").Value = 156
$ws.Range("H14").Value = 18132

# Oman (row 35)
$ws.Range("B35").Value = 81067
$ws.Range("C35").Value = 354
$ws.Range("D35").Value = 72263
$ws.Range("E35").Value = 8302

# Israel (row 36)
$ws.Range("B36").Value = 80431
$ws.Range("C36").Value = 872
$ws.Range("D36").Value = 55274
$ws.Range("E36").Value = 24579
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 578

# Rumania (row 44)
$ws.Range("B44").Value = 59273
$ws.Range("C44").Value = 1378
$ws.Range("D44").Value = 29289
$ws.Range("E44").Value = 27368
$ws.Range("G44").Value = 50
$ws.Range("H44").Value = 2616

# Suiza (row 58)
$ws.Range("B58").Value = 36269
$ws.Range("C58").Value = 161
$ws.Range("E58").Value = 2683
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1986

# El Salvador (row 73)
$ws.Range("B73").Value = 19544
$ws.Range("C73").Value = 418
$ws.Range("D73").Value = 9379
$ws.Range("E73").Value = 9645

# Malasia (row 88)
$ws.Range("B88").Value = 9063
$ws.Range("C88").Value = 25
$ws.Range("D88").Value = 8728
$ws.Range("E88").Value = 210

# Row 151: Jamaica -> Malta (Malta's refreshed stats move it above Jamaica)
$ws.Range("A151").Value = "Malta"
$ws.Range("B151").Value = 995
$ws.Range("C151").Value = 49
$ws.Range("D151").Value = 675
$ws.Range("E151").Value = 311
$ws.Range("H151").Value = 9

# Row 152: Malta -> Jamaica (Jamaica's stats are unchanged, just re-ranked)
$ws.Range("A152").Value = "Jamaica"
$ws.Range("B152").Value = 958
$ws.Range("C152").Value = 30
$ws.Range("D152").Value = 745
$ws.Range("E152").Value = 201
$ws.Range("H152").Value = 12

# Curazao (row 200)
$ws.Range("D200").Value = 30
$ws.Range("E200").Value = 0

# Row 202/203: Santa Lucia <-> Timor Oriental swap places (tied stats)
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"
